$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu values for Case_3_191 (380 kV case)

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.053645271017139
$ws.Cells.Item(2, 4).Value = 1.051910786139508
$ws.Cells.Item(2, 5).Value = 1.059894502419626
$ws.Cells.Item(2, 6).Value = 1.070230638277563
$ws.Cells.Item(2, 9).Value = 1.040300881744945
$ws.Cells.Item(2, 10).Value = 1.058661527357204
$ws.Cells.Item(2, 11).Value = 1.054660978389271
$ws.Cells.Item(2, 12).Value = 1.062622762400373
$ws.Cells.Item(2, 13).Value = 1.072931029374531
$ws.Cells.Item(2, 14).Value = 1.060164948156226

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.055498791462641
$ws.Cells.Item(3, 4).Value = 1.053332612497063
$ws.Cells.Item(3, 5).Value = 1.061634772963499
$ws.Cells.Item(3, 6).Value = 1.072230342341859
$ws.Cells.Item(3, 9).Value = 1.040794297841773
$ws.Cells.Item(3, 10).Value = 1.060161631065802
$ws.Cells.Item(3, 11).Value = 1.05589398107142
$ws.Cells.Item(3, 12).Value = 1.064175013051718
$ws.Cells.Item(3, 13).Value = 1.074744129012214
$ws.Cells.Item(3, 14).Value = 1.061667182184154

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.056694124493594
$ws.Cells.Item(4, 4).Value = 1.054248992432403
$ws.Cells.Item(4, 5).Value = 1.062756893120052
$ws.Cells.Item(4, 6).Value = 1.073520551541594
$ws.Cells.Item(4, 9).Value = 1.0411102614864
$ws.Cells.Item(4, 10).Value = 1.061128037080686
$ws.Cells.Item(4, 11).Value = 1.056687666719927
$ws.Cells.Item(4, 12).Value = 1.065175023784708
$ws.Cells.Item(4, 13).Value = 1.075913192450407
$ws.Cells.Item(4, 14).Value = 1.062634960606428

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.057195701103482
$ws.Cells.Item(5, 4).Value = 1.054633382637259
$ws.Cells.Item(5, 5).Value = 1.063227705943902
$ws.Cells.Item(5, 6).Value = 1.074062083982648
$ws.Cells.Item(5, 9).Value = 1.041242306730834
$ws.Cells.Item(5, 10).Value = 1.061533311047924
$ws.Cells.Item(5, 11).Value = 1.057020352657722
$ws.Cells.Item(5, 12).Value = 1.065594393312044
$ws.Cells.Item(5, 13).Value = 1.076403697882963
$ws.Cells.Item(5, 14).Value = 1.063040810109185

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.057279863415789
$ws.Cells.Item(6, 4).Value = 1.05469787367735
$ws.Cells.Item(6, 5).Value = 1.063306703720013
$ws.Cells.Item(6, 6).Value = 1.074152959197178
$ws.Cells.Item(6, 9).Value = 1.041264431810147
$ws.Cells.Item(6, 10).Value = 1.061601299998955
$ws.Cells.Item(6, 11).Value = 1.057076155052369
$ws.Cells.Item(6, 12).Value = 1.065664747134078
$ws.Cells.Item(6, 13).Value = 1.076485999686909
$ws.Cells.Item(6, 14).Value = 1.063108895612324

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.056700830255807
$ws.Cells.Item(7, 4).Value = 1.054254132016255
$ws.Cells.Item(7, 5).Value = 1.062763187758131
$ws.Cells.Item(7, 6).Value = 1.073527790918917
$ws.Cells.Item(7, 9).Value = 1.04111202896056
$ws.Cells.Item(7, 10).Value = 1.061133456289677
$ws.Cells.Item(7, 11).Value = 1.056692115915603
$ws.Cells.Item(7, 12).Value = 1.065180631462944
$ws.Cells.Item(7, 13).Value = 1.075919750386834
$ws.Cells.Item(7, 14).Value = 1.062640387511318

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.054272521525432
$ws.Cells.Item(8, 4).Value = 1.05239206095627
$ws.Cells.Item(8, 5).Value = 1.060483464898053
$ws.Cells.Item(8, 6).Value = 1.070907234006684
$ws.Cells.Item(8, 9).Value = 1.040468323195654
$ws.Cells.Item(8, 10).Value = 1.059169387106455
$ws.Cells.Item(8, 11).Value = 1.055078545021187
$ws.Cells.Item(8, 12).Value = 1.063148274678977
$ws.Cells.Item(8, 13).Value = 1.073544643215066
$ws.Cells.Item(8, 14).Value = 1.060673529124573

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.04996181931677
$ws.Cells.Item(9, 4).Value = 1.049082333873399
$ws.Cells.Item(9, 5).Value = 1.056435161647647
$ws.Cells.Item(9, 6).Value = 1.066259885126166
$ws.Cells.Item(9, 9).Value = 1.039308392137335
$ws.Cells.Item(9, 10).Value = 1.055675037774338
$ws.Cells.Item(9, 11).Value = 1.052202817781133
$ws.Cells.Item(9, 12).Value = 1.059532487297663
$ws.Cells.Item(9, 13).Value = 1.069326811291963
$ws.Cells.Item(9, 14).Value = 1.057174217415596

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.047065382345358
$ws.Cells.Item(10, 4).Value = 1.046855734634067
$ws.Cells.Item(10, 5).Value = 1.053714120154176
$ws.Cells.Item(10, 6).Value = 1.063140346875264
$ws.Cells.Item(10, 9).Value = 1.038517454039062
$ws.Cells.Item(10, 10).Value = 1.053321923618442
$ws.Cells.Item(10, 11).Value = 1.05026299010934
$ws.Cells.Item(10, 12).Value = 1.057097599654044
$ws.Cells.Item(10, 13).Value = 1.066491698632289
$ws.Cells.Item(10, 14).Value = 1.054817761567693

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.045805529561924
$ws.Cells.Item(11, 4).Value = 1.045886610542765
$ws.Cells.Item(11, 5).Value = 1.052530343005545
$ws.Cells.Item(11, 6).Value = 1.061784183135766
$ws.Cells.Item(11, 9).Value = 1.038170689634268
$ws.Cells.Item(11, 10).Value = 1.052297170166802
$ws.Cells.Item(11, 11).Value = 1.04941744562471
$ws.Cells.Item(11, 12).Value = 1.056037228474749
$ws.Cells.Item(11, 13).Value = 1.065258260545546
$ws.Cells.Item(11, 14).Value = 1.053791552848612

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.045336685762857
$ws.Cells.Item(12, 4).Value = 1.045525866369943
$ws.Cells.Item(12, 5).Value = 1.052089777884998
$ws.Cells.Item(12, 6).Value = 1.061279607434165
$ws.Cells.Item(12, 9).Value = 1.038041234302408
$ws.Cells.Item(12, 10).Value = 1.051915632624488
$ws.Cells.Item(12, 11).Value = 1.049102516152185
$ws.Cells.Item(12, 12).Value = 1.055642428144239
$ws.Cells.Item(12, 13).Value = 1.064799207476943
$ws.Cells.Item(12, 14).Value = 1.053409473479224

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.045437294511967
$ws.Cells.Item(13, 4).Value = 1.0456032822859
$ws.Cells.Item(13, 5).Value = 1.052184319795392
$ws.Cells.Item(13, 6).Value = 1.061387878913749
$ws.Cells.Item(13, 9).Value = 1.038069032528171
$ws.Cells.Item(13, 10).Value = 1.051997514723741
$ws.Cells.Item(13, 11).Value = 1.049170108637024
$ws.Cells.Item(13, 12).Value = 1.055727156661471
$ws.Cells.Item(13, 13).Value = 1.064897717024573
$ws.Cells.Item(13, 14).Value = 1.05349147186045

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.045766792855075
$ws.Cells.Item(14, 4).Value = 1.045856807097862
$ws.Cells.Item(14, 5).Value = 1.052493943389786
$ws.Cells.Item(14, 6).Value = 1.061742491966795
$ws.Cells.Item(14, 9).Value = 1.038160002155875
$ws.Cells.Item(14, 10).Value = 1.052265650606855
$ws.Cells.Item(14, 11).Value = 1.049391431031856
$ws.Cells.Item(14, 12).Value = 1.056004613285039
$ws.Cells.Item(14, 13).Value = 1.065220333568899
$ws.Cells.Item(14, 14).Value = 1.053759988527273

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.045969690494384
$ws.Cells.Item(15, 4).Value = 1.04601290973628
$ws.Cells.Item(15, 5).Value = 1.052684598343098
$ws.Cells.Item(15, 6).Value = 1.061960869155429
$ws.Cells.Item(15, 9).Value = 1.038215964958448
$ws.Cells.Item(15, 10).Value = 1.052430738277877
$ws.Cells.Item(15, 11).Value = 1.049527681032961
$ws.Cells.Item(15, 12).Value = 1.056175439388219
$ws.Cells.Item(15, 13).Value = 1.065418988341461
$ws.Cells.Item(15, 14).Value = 1.053925310641725

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.047148872855715
$ws.Cells.Item(16, 4).Value = 1.046919945473602
$ws.Cells.Item(16, 5).Value = 1.053792564607353
$ws.Cells.Item(16, 6).Value = 1.06323023525426
$ws.Cells.Item(16, 9).Value = 1.038540376722552
$ws.Cells.Item(16, 10).Value = 1.053389808298988
$ws.Cells.Item(16, 11).Value = 1.050318986942183
$ws.Cells.Item(16, 12).Value = 1.057167843642872
$ws.Cells.Item(16, 13).Value = 1.066573433158028
$ws.Cells.Item(16, 14).Value = 1.054885742652273

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.047887005118885
$ws.Cells.Item(17, 4).Value = 1.047487555674341
$ws.Cells.Item(17, 5).Value = 1.054486060301619
$ws.Cells.Item(17, 6).Value = 1.06402501467124
$ws.Cells.Item(17, 9).Value = 1.038742719279361
$ws.Cells.Item(17, 10).Value = 1.053989829339624
$ws.Cells.Item(17, 11).Value = 1.050813844471042
$ws.Cells.Item(17, 12).Value = 1.057788717399388
$ws.Cells.Item(17, 13).Value = 1.067296010988358
$ws.Cells.Item(17, 14).Value = 1.055486615791609

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.048316999227654
$ws.Cells.Item(18, 4).Value = 1.047818152706113
$ws.Cells.Item(18, 5).Value = 1.054890030922801
$ws.Cells.Item(18, 6).Value = 1.064488078534763
$ws.Cells.Item(18, 9).Value = 1.038860329578164
$ws.Cells.Item(18, 10).Value = 1.054339249479351
$ws.Cells.Item(18, 11).Value = 1.051101948681325
$ws.Cells.Item(18, 12).Value = 1.058150280474943
$ws.Cells.Item(18, 13).Value = 1.067716919032015
$ws.Cells.Item(18, 14).Value = 1.055836532148014

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.048463524271725
$ws.Cells.Item(19, 4).Value = 1.04793079681849
$ws.Cells.Item(19, 5).Value = 1.055027684664918
$ws.Cells.Item(19, 6).Value = 1.064645884386819
$ws.Cells.Item(19, 9).Value = 1.038900361894941
$ws.Cells.Item(19, 10).Value = 1.054458298031238
$ws.Cells.Item(19, 11).Value = 1.051200094069514
$ws.Cells.Item(19, 12).Value = 1.058273466085122
$ws.Cells.Item(19, 13).Value = 1.067860343684066
$ws.Cells.Item(19, 14).Value = 1.055955749762499

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.047807867167125
$ws.Cells.Item(20, 4).Value = 1.047426706293342
$ws.Cells.Item(20, 5).Value = 1.054411710151308
$ws.Cells.Item(20, 6).Value = 1.063939796088486
$ws.Cells.Item(20, 9).Value = 1.038721052583045
$ws.Cells.Item(20, 10).Value = 1.053925511041361
$ws.Cells.Item(20, 11).Value = 1.050760806694881
$ws.Cells.Item(20, 12).Value = 1.057722163890803
$ws.Cells.Item(20, 13).Value = 1.067218543268643
$ws.Cells.Item(20, 14).Value = 1.055422206153986

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.045669788260189
$ws.Cells.Item(21, 4).Value = 1.04578217171823
$ws.Cells.Item(21, 5).Value = 1.052402790826381
$ws.Cells.Item(21, 6).Value = 1.06163809056336
$ws.Cells.Item(21, 9).Value = 1.038133231920854
$ws.Cells.Item(21, 10).Value = 1.052186716246404
$ws.Cells.Item(21, 11).Value = 1.049326280889729
$ws.Cells.Item(21, 12).Value = 1.05592293511448
$ws.Cells.Item(21, 13).Value = 1.065125356067829
$ws.Cells.Item(21, 14).Value = 1.053680942070977

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.044320398231056
$ws.Cells.Item(22, 4).Value = 1.044743732921152
$ws.Cells.Item(22, 5).Value = 1.0511347296775
$ws.Cells.Item(22, 6).Value = 1.060186065668835
$ws.Cells.Item(22, 9).Value = 1.037759871713805
$ws.Cells.Item(22, 10).Value = 1.051088257894297
$ws.Cells.Item(22, 11).Value = 1.048419372770887
$ws.Cells.Item(22, 12).Value = 1.054786289203849
$ws.Cells.Item(22, 13).Value = 1.063804070547705
$ws.Cells.Item(22, 14).Value = 1.052580923782014

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.045036226986692
$ws.Cells.Item(23, 4).Value = 1.045294657421848
$ws.Cells.Item(23, 5).Value = 1.05180743245442
$ws.Cells.Item(23, 6).Value = 1.06095628061607
$ws.Cells.Item(23, 9).Value = 1.037958157502014
$ws.Cells.Item(23, 10).Value = 1.051671072353573
$ws.Cells.Item(23, 11).Value = 1.048900618381046
$ws.Cells.Item(23, 12).Value = 1.055389366099847
$ws.Cells.Item(23, 13).Value = 1.064505011996365
$ws.Cells.Item(23, 14).Value = 1.053164565904673

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.047843627887323
$ws.Cells.Item(24, 4).Value = 1.047454202993353
$ws.Cells.Item(24, 5).Value = 1.054445307432262
$ws.Cells.Item(24, 6).Value = 1.063978304296186
$ws.Cells.Item(24, 9).Value = 1.038730844105478
$ws.Cells.Item(24, 10).Value = 1.053954575447722
$ws.Cells.Item(24, 11).Value = 1.050784773844325
$ws.Cells.Item(24, 12).Value = 1.057752238354538
$ws.Cells.Item(24, 13).Value = 1.06725354932435
$ws.Cells.Item(24, 14).Value = 1.055451311835137

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.051080131415635
$ws.Cells.Item(25, 4).Value = 1.049941453308942
$ws.Cells.Item(25, 5).Value = 1.057485561581816
$ws.Cells.Item(25, 6).Value = 1.067464987057414
$ws.Cells.Item(25, 9).Value = 1.039611342568615
$ws.Cells.Item(25, 10).Value = 1.056582480805428
$ws.Cells.Item(25, 11).Value = 1.052950193307495
$ws.Cells.Item(25, 12).Value = 1.060471463370342
$ws.Cells.Item(25, 13).Value = 1.070421218135995
$ws.Cells.Item(25, 14).Value = 1.058082949119876

Write-Output "vm_pu values updated for 380 kV case"
